# 16.3.1.xlsx — add a 2023 data column (E) next to the existing 2018 column (D),
# and update the source citation (row 9) to reference both survey years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "2023" column -------------------------------------------------
# Year header, alongside the existing 2018 header in D4.
$ws.Range("E4").Value = 2023

# Overall indicator value for 2023 (alongside D5 = 29.5 for 2018).
$ws.Range("E5").Value = 38

# "Жерлери / Местность / Urbanisation" sub-header row: no 2023 total, but
# make the whole row bold (matches the new bold styling of A6:E6). Clone
# D6's formatting onto the new E6 cell first so the bolding below produces
# the same font as D6 instead of falling back to the default Calibri font.
$ws.Range("D6").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E6").Value = ""
$ws.Range("A6:E6").Font.Bold = $true

# Urban ("шаар" / "Городская" / "Urban") and Rural ("айыл" / "Сельская" /
# "Rural") rows: no numeric breakdown is available for 2023, so the new
# cells just hold a dash, right aligned like the rest of the data column.
# Clone D7/D8's formatting (plain Times New Roman 9pt) onto E7/E8 so the
# font matches instead of defaulting to Calibri.
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E7").Value = "-"
$ws.Range("E7").HorizontalAlignment = -4152   # xlRight

$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E8").Value = "-"
$ws.Range("E8").HorizontalAlignment = -4152   # xlRight

# --- Update the footnote / source citation row (row 9) -----------------
$ws.Range("A9").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2018-ж., 2023-ж."
$ws.Range("B9").Value = "По данным кластерного обследования по многим показателям, 2018г., 2023г."
$ws.Range("C9").Value = "According to the cluster survey in many respects, 2018, 2023."
